$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full contents of row 13 and row 14 (these two
# observation records were reordered). Only the columns whose values
# actually differ between the two rows are listed here; the rest of the
# row (D, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY) already
# holds identical values on both rows, so nothing else needs touching.

$cols = @("A","B","E","F","G","H","Q","R","X","Z","AB")

foreach ($col in $cols) {
    $cell13 = $ws.Range("$col`13")
    $cell14 = $ws.Range("$col`14")
    $v13 = $cell13.Value2
    $v14 = $cell14.Value2
    $cell13.Value = $v14
    $cell14.Value = $v13
}
